$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of order data to append after existing row 11
$newRows = @(
    @("AH252", "Natalie's - Orange Juice", "1", "24.50", "24.50"),
    @("TN454", "Natalie's - Orange Mango", "1", "13.38", "13.38"),
    @("TN330", "Natalie's - Honey Tangerine", "1", "14.57", "14.57")
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $rowRange = $ws.Range("A" + $r + ":E" + $r)
    # Force text storage (matches source data which is typed as text/inlineStr),
    # then restore the default "Normal" style so no new style index is introduced.
    $rowRange.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $rowRange.Style = "Normal"
}
